$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- New row of data (row 2) --
# A2: date/time serial value, formatted as a date (maps to builtin numFmtId 22)
$ws.Range("A2").Value = 42587.832407407404

# B2: new shared string "Noun" (Method column)
$ws.Range("B2").Value = "Noun"

# C2:M2: plain numeric stats
$ws.Range("C2").Value = 3161
$ws.Range("D2").Value = 75
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 66
$ws.Range("I2").Value = 33
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 57
$ws.Range("M2").Value = 42

# Apply a date/time number format to column A (header + new data cell) -
# built-in format code 22 ("m/d/yy h:mm")
$ws.Range("A1:A2").NumberFormat = "m/d/yy h:mm"

# Widen column A so the date values are fully visible
$ws.Columns("A").ColumnWidth = 13
